$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '80.882.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.40%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.133.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.81%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.05'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '614.53'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.90%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.277'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +22.62%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.574'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.30%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.132.56'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.82%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.569'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.57%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000248'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +10.96%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.00%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.713.00'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.83%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.92%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.838.90'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.50%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.137.99'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.51%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.11'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +8.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.75'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.87%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.61'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.87'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.64%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.03'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.12%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.81%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.10'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +7.97%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.293.29'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.77%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '75.42'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.72'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.50%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.31%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000119'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.41%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.18%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.83'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.46%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '549.04'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.99%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.42%  '

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.62%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.21%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.38'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.05%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.08%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.91'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +9.56%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.66%  '

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.99'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +12.34%  '

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.98'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +19.50%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '158.91'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.18%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '186.05'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.28%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.70'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.72%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.23%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.77%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.43'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.90%  '
